$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.564.91"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "1.867.37"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'324.36"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.4607"
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("D8").Value = "'0.3868"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.07854"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "'0.9736"
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "1.890.81"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "'6.972"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "'5.691"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "'0.06934"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").Value = "'88.11"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "28.565.95"
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("D22").Value = "'5.262"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("D23").Value = "'11.02"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "'2.110"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "2.073.65"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'152.37"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").Value = "'19.22"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "'5.767"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").Value = "'1.983"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").Value = "'119.16"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").Value = "'0.09328"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").Value = "'0.9175"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "'5.257"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").Value = "'1.332"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("D35").Value = "'3.326"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").Value = "'0.05783"
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("D37").Value = "'1.154"
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").Value = "'0.02075"
$ws.Range("E38").Value = "  -3.17%  "
$ws.Range("D39").Value = "'7.726"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "'0.5619"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").Value = "'0.1782"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "'9.766"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "'0.07170"
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("D44").Value = "'11.76"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").Value = "'0.5291"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "'2.141"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "'1.133"
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("D48").Value = "'1.831"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").Value = "'112.73"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").Value = "'2.414"
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("E51").Value = "  +0.12%  "
